# Append two new geocoding-result rows (14 and 15) to Sheet1, matching the
# pattern of the existing rows: Latitude/Longitude (A/B) blank, Full Address
# (C) filled in -- same "Karmayogi Bhavan..." address already seen in the
# sheet (rows 9 and 13).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$addr = "Karmayogi Bhavan, Road 3, Sector 7, Gandhinagar, Gandhinagar Taluka, Gandhinagar, Gujarat, 382008, India"

# A leading single-quote forces an empty *text* cell (instead of Value=""
# which just clears/blanks the cell outright) so A/B end up the same
# empty-string-typed cells as the rest of the column; resetting the Style
# back to "Normal" drops the quote-prefix formatting Excel would otherwise
# apply, so the cell is plain text with no visible marker.
$ws.Range("A14").Value = "'"
$ws.Range("A14").Style = "Normal"
$ws.Range("B14").Value = "'"
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Value = $addr

$ws.Range("A15").Value = "'"
$ws.Range("A15").Style = "Normal"
$ws.Range("B15").Value = "'"
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Value = $addr
